# Update "想去人数"(F) / "最低票价"(G) figures across all four sheets
# to match the refreshed scrape output (gh-pages data refresh at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 805
$ws.Range("F3").Value = 976
$ws.Range("F4").Value = 760
$ws.Range("F5").Value = 861
$ws.Range("F6").Value = 424
$ws.Range("F7").Value = 657
$ws.Range("F9").Value = 1248
$ws.Range("F10").Value = 677
$ws.Range("F11").Value = 402
$ws.Range("F12").Value = 533
$ws.Range("F13").Value = 175
$ws.Range("F14").Value = 23
$ws.Range("F15").Value = 765
$ws.Range("F17").Value = 383
$ws.Range("F18").Value = 365
$ws.Range("F19").Value = 87
$ws.Range("F20").Value = 569
$ws.Range("F21").Value = 123
$ws.Range("F22").Value = 614
$ws.Range("F23").Value = 32
$ws.Range("F24").Value = 885

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = "不可售"
$ws.Range("F4").Value = 109
$ws.Range("F7").Value = 186
$ws.Range("F10").Value = 28
$ws.Range("F12").Value = 104

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 372

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 372
$ws.Range("G3").Value = "不可售"
$ws.Range("F5").Value = 805
$ws.Range("F6").Value = 976
$ws.Range("F7").Value = 760
$ws.Range("F8").Value = 861
$ws.Range("F9").Value = 424
$ws.Range("F10").Value = 657
$ws.Range("F12").Value = 1248
$ws.Range("F13").Value = 677
$ws.Range("F14").Value = 109
$ws.Range("F16").Value = 402
$ws.Range("F17").Value = 533
$ws.Range("F19").Value = 175
$ws.Range("F20").Value = 23
$ws.Range("F21").Value = 765
$ws.Range("F22").Value = 186
$ws.Range("F24").Value = 383
$ws.Range("F25").Value = 365
$ws.Range("F26").Value = 87
$ws.Range("F29").Value = 569
$ws.Range("F30").Value = 28
$ws.Range("F32").Value = 104
$ws.Range("F33").Value = 104
$ws.Range("F34").Value = 123
$ws.Range("F35").Value = 614
$ws.Range("F36").Value = 32
$ws.Range("F37").Value = 885
